$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Important Notes" text in A1: change "65" to "55" hours for middle school PD
$ws.Range("A1").Value = "Important Notes: `n1) Code.org will match up to 1/2 of the stipends costs for the district. `n2) Code.org will not cover teacher taxes and benfits. `n3) Code.org will pay the district directly, which needs to, in turn, pay the teachers for their time in professional development (PD). `n4) This calcuator assumes that all teachers will attend all phases of the PD. See the application for detailed description of the PD model.`n5) Total hours for high school PD is 94 (online and in-person) and total hours for middle school PD is 55 (online and in-person).`n6) The number of committed teachers entered below should be the TOTAL number of teachers you plan to enroll in either high school course (Exploring Computer Science or Computer Science Principles) or either middle school courses (Computer Science in Algebra, Computer Science in Science). For example, you may want to enroll 30 ECS teachers and 30 CSP teachers therefore you would enter 60 teachers in cell B4."

# Update the formulas for Middle School in C8 and C9: 65 -> 55
$ws.Range("C8").Formula = "=(55*B3*B5)/2"
$ws.Range("C9").Formula = "=(55*B3*B5)/2"

# Update the selected cell to C10
$ws.Range("C10").Select()
